$d = $word.ActiveDocument

# --- Create the three new character styles (in the order they appear in the diff) ---

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.NameAscii = "Calibri"
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.NameAscii = "Calibri"
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.NameAscii = "Calibri"
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas das campanhas ..." run (4 occurrences) ---

$datasText = "Datas das campanhas de 2022 que usam Constelação de Cygnus: 10 a 19 de agosto, 9 a 18 de setembro, 8 a 17 de outubro"
$rng = $d.Content
while ($rng.Find.Execute($datasText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "Está a participar numa campanha global ..." run ---

$estaText = "Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  Constelação de Cygnus no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo."
$rng = $d.Content
if ($rng.Find.Execute($estaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "por Jenik Hollan, ..." run ---

$jenikText = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng = $d.Content
if ($rng.Find.Execute($jenikText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
